$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: replace the "nombre" / "edad" header pair with a single free-text cell in A1.
$ws.Range("A1").Value = "Aqui tengo una mini tabla para que me la hagas excel:"

# Row 2: A2 becomes "edad" (was "Lucia").
$ws.Range("A2").Value = "edad"

# Row 3: A3 becomes the text "25" (was "Pedro"); keep it as text, not a number.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "25"
$ws.Range("A3").ClearFormats()

# Row 4 (new): A4 becomes the text "33".
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "33"
$ws.Range("A4").ClearFormats()

# Column B is no longer used anywhere (B1 "edad", B2 25, B3 33 all removed).
$ws.Range("B1:B3").ClearFormats()
$ws.Range("B1:B3").Value = $null
